$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B30 value changed from 5 to 0
$ws.Range("B30").Value = 0

# A new channel row (7033) was appended as row 33, mirroring the existing
# "channel code" rows above it (text value in column A, count 0 in column B).
# Column A holds these codes as text (shared strings) even though they look
# numeric, so we build the text value via a formula in a scratch cell and
# paste only the value into A33 - this avoids Excel auto-converting the
# literal "7033" into a numeric cell.
$ws.Range("D1").Formula = "=""7033"""
$ws.Range("D1").Copy()
$ws.Range("A33").PasteSpecial(-4163)  # xlPasteValues
$ws.Range("D1").ClearContents()

$ws.Range("B33").Value = 0
